$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: remove Trial 2..Trial 10 (columns C..K), shift Mean/Stdev to C/D
$ws.Range("C1").Value = "Mean"
$ws.Range("D1").Value = "Stdev"
$ws.Range("E1:M1").Clear()

# Row 2: Window Size:20 (overwrites the old A2 "Window Size:1" text, which moves to A10 below)
$ws.Range("A2").Value = "SlidingLin with Window Size:20"
$ws.Range("B2").Value = 85.22132138585457
$ws.Range("C2").Formula = "=AVERAGE(B2:B2)"
$ws.Range("D2").Formula = "=STDEV(B2:B2)"

# Row 3: Window Size:15
$ws.Range("A3").Value = "SlidingLin with Window Size:15"
$ws.Range("B3").Value = 113.4426791493449
$ws.Range("C3").Formula = "=AVERAGE(B3:B3)"
$ws.Range("D3").Formula = "=STDEV(B3:B3)"

# Row 4: Window Size:10
$ws.Range("A4").Value = "SlidingLin with Window Size:10"
$ws.Range("B4").Value = 4020.192172002387
$ws.Range("C4").Formula = "=AVERAGE(B4:B4)"
$ws.Range("D4").Formula = "=STDEV(B4:B4)"

# Row 5 intentionally left blank (gap in diff)

# Row 6: Window Size:5
$ws.Range("A6").Value = "SlidingLin with Window Size:5"
$ws.Range("B6").Value = 821473477.520296
$ws.Range("C6").Formula = "=AVERAGE(B6:B6)"
$ws.Range("D6").Formula = "=STDEV(B6:B6)"

# Row 7: Window Size:4
$ws.Range("A7").Value = "SlidingLin with Window Size:4"
$ws.Range("B7").Value = 133.6683580942512
$ws.Range("C7").Formula = "=AVERAGE(B7:B7)"
$ws.Range("D7").Formula = "=STDEV(B7:B7)"

# Row 8: Window Size:3 (text only)
$ws.Range("A8").Value = "SlidingLin with Window Size:3"

# Row 9: Window Size:2 (text only)
$ws.Range("A9").Value = "SlidingLin with Window Size:2"

# Row 10: Window Size:1 (text only) - moved from A2
$ws.Range("A10").Value = "SlidingLin with Window Size:1"
